$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-TextValue "D2" "20.410.06"
$ws.Range("E2").Value = "  -6.33%  "
Set-TextValue "D3" "1.438.41"
$ws.Range("E3").Value = "  -6.70%  "
$ws.Range("E4").Value = "  -0.66%  "
Set-TextValue "D5" "1.004"
$ws.Range("E5").Value = "  -0.18%  "
Set-TextValue "D6" "277.64"
$ws.Range("E6").Value = "  -3.81%  "
Set-TextValue "D7" "0.3745"
$ws.Range("E7").Value = "  -3.65%  "
Set-TextValue "D8" "0.3076"
$ws.Range("E8").Value = "  -3.51%  "
Set-TextValue "D9" "40.39"
$ws.Range("E9").Value = "  -6.86%  "
$ws.Range("E10").Value = "  -4.23%  "
Set-TextValue "D11" "0.06575"
$ws.Range("E11").Value = "  -7.86%  "
Set-TextValue "D12" "1.002"
$ws.Range("E12").Value = "  -0.68%  "
Set-TextValue "D13" "5.351"
$ws.Range("E13").Value = "  -4.07%  "
Set-TextValue "D14" "17.30"
$ws.Range("E14").Value = "  -6.05%  "
Set-TextValue "D15" "6.124"
$ws.Range("E15").Value = "  -7.32%  "
Set-TextValue "D16" "1.440.34"
$ws.Range("E16").Value = "  -6.63%  "
Set-TextValue "D17" "0.00001011"
$ws.Range("E17").Value = "  -7.33%  "
Set-TextValue "D18" "76.92"
$ws.Range("E18").Value = "  -7.29%  "
Set-TextValue "D19" "0.05821"
$ws.Range("E19").Value = "  -11.08%  "
Set-TextValue "D20" "1.003"
$ws.Range("E20").Value = "  -0.01%  "
Set-TextValue "D21" "5.716"
$ws.Range("E21").Value = "  -6.45%  "
Set-TextValue "D22" "14.37"
$ws.Range("E22").Value = "  -5.29%  "
Set-TextValue "D23" "10.81"
$ws.Range("E23").Value = "  -1.16%  "
Set-TextValue "D24" "2.324"
$ws.Range("E24").Value = "  -2.15%  "
Set-TextValue "D25" "20.404.97"
$ws.Range("E25").Value = "  -6.48%  "
Set-TextValue "D26" "142.32"
$ws.Range("E26").Value = "  -1.91%  "
Set-TextValue "D27" "2.211"
$ws.Range("E27").Value = "  -6.13%  "
Set-TextValue "D28" "17.02"
$ws.Range("E28").Value = "  -7.43%  "
Set-TextValue "D29" "1.602.17"
$ws.Range("E29").Value = "  -6.80%  "
Set-TextValue "D30" "109.92"
$ws.Range("E30").Value = "  -6.01%  "
Set-TextValue "D31" "3.930"
$ws.Range("E31").Value = "  -19.10%  "
Set-TextValue "D32" "0.9096"
$ws.Range("E32").Value = "  -6.06%  "
Set-TextValue "D33" "5.457"
$ws.Range("E33").Value = "  -6.56%  "
Set-TextValue "D34" "0.07702"
$ws.Range("E34").Value = "  -5.64%  "
Set-TextValue "D35" "8.361"
$ws.Range("E35").Value = "  -6.49%  "
Set-TextValue "D36" "1.002"
$ws.Range("E36").Value = "  -0.01%  "
Set-TextValue "D37" "10.87"
$ws.Range("E37").Value = "  +3.62%  "
Set-TextValue "D38" "0.05669"
$ws.Range("E38").Value = "  -5.38%  "
$ws.Range("E39").Value = "  -2.93%  "
Set-TextValue "D40" "4.719"
$ws.Range("E40").Value = "  -6.60%  "
$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D41" "0.1913"
$ws.Range("E41").Value = "  -5.35%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D42" "0.02034"
$ws.Range("E42").Value = "  -8.36%  "
Set-TextValue "D43" "1.313"
$ws.Range("E43").Value = "  -17.53%  "
Set-TextValue "D44" "3.585"
$ws.Range("E44").Value = "  -3.71%  "
Set-TextValue "D45" "0.5319"
$ws.Range("E45").Value = "  -6.61%  "
Set-TextValue "D46" "12.00"
$ws.Range("E46").Value = "  -6.75%  "
Set-TextValue "D47" "0.5148"
$ws.Range("E47").Value = "  -6.27%  "
Set-TextValue "D48" "111.95"
$ws.Range("E48").Value = "  -3.53%  "
Set-TextValue "D49" "1.789"
$ws.Range("E49").Value = "  -3.32%  "
Set-TextValue "D50" "1.055"
$ws.Range("E50").Value = "  -6.16%  "
$ws.Range("E51").Value = "  +0.07%  "
